$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 9): #property=8, address=<same as row8>, date=10/1/2024,
# numberOfPrimaryInsured=120114, numberOfPeopleCertifiedForLongTermCare=25994
$ws.Range("A9").Value = 8
$ws.Range("D9").Value = 120114
$ws.Range("E9").Value = 25994

# Copy the formatting of row 8 onto row 9 so the new cells pick up the same
# cell styles (A/D/E use style index 1, B/C stay unstyled) without minting
# any new style entries in styles.xml.
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)  # xlPasteFormats

# B9 repeats the same address string already used in B8 (shared string).
$ws.Range("B9").Value = $ws.Range("B8").Text

# C9 needs the literal text "10/1/2024" (not an auto-converted date serial).
# Stage it in a scratch cell via a formula (so no date-parsing happens and no
# number-format style gets created), copy its computed value into C9, then
# remove the scratch column again.
$ws.Range("G9").Formula = "=""10/1/2024"""
$ws.Range("G9").Copy()
$ws.Range("C9").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("G9").EntireColumn.Delete()

$ws.Range("F10").Select()
